$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = "Fri Mar 08 01:09:01 EST 2024"
$ws.Range("B3").Value = "Fri Mar 08 01:09:40 EST 2024"
$ws.Range("B5").Value = "Fri Mar 08 01:10:18 EST 2024"
$ws.Range("B6").Value = "Fri Mar 08 01:11:30 EST 2024"
$ws.Range("B7").Value = "Fri Mar 08 01:12:11 EST 2024"
